$d = $word.ActiveDocument

# The template shipped a pile of near-duplicate "egXML" / "egXMLTable"
# paragraph styles (egXML, egXML0..egXML4, egXMLTable, egXMLTable0..egXMLTable4)
# - leftovers from repeated imports of the same TEI gloss-list styling.
# Clean them all up, keeping only a single canonical "egXML" and
# "egXMLTable" style pair (re-created fresh so the surviving style carries
# the plain "egXML" / "egXMLTable" id instead of one of the numbered
# duplicates).

$duplicateStyleIds = @(
    "egXMLTable",
    "egXML",
    "egXML0",
    "egXMLTable0",
    "egXML1",
    "egXMLTable1",
    "egXML2",
    "egXMLTable2",
    "egXML3",
    "egXMLTable3",
    "egXML4",
    "egXMLTable4"
)

foreach ($styleId in $duplicateStyleIds) {
    $existing = $d.Styles($styleId)
    if ($existing -ne $null) {
        $existing.Delete()
    }
}

# Re-create the single canonical "egXML" paragraph style.
$egXml = $d.Styles.Add("egXML", 1)
$egXml.BaseStyle = $d.Styles("Normal")
$egXml.QuickStyle = $true
$egXml.Font.Name = "Courier"
$egXml.Font.Size = 10

# Re-create the single canonical "egXMLTable" paragraph style.
$egXmlTable = $d.Styles.Add("egXMLTable", 1)
$egXmlTable.BaseStyle = $d.Styles("Normal")
$egXmlTable.QuickStyle = $true
$egXmlTable.Font.Name = "Courier"
$egXmlTable.Font.Size = 9
$egXmlTable.ParagraphFormat.SpaceBefore = 4
